$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$headers = @("pairs", "Df", "SumsOfSqs", "F.Model", "R2", "p.value", "p.adjusted", "sig")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Data rows: pairs, Df, SumsOfSqs, F.Model, R2, p.value, p.adjusted, sig
$data = @(
    @("pre_ltx vs healthy", 1, 1222.92106840471, 7.94994650341597, 0.0187887705818079, 0.001, 0.0015, "**"),
    @("pre_ltx vs post_ltx", 1, 630.146146437502, 4.06180919840911, 0.00662047049736794, 0.008, 0.008, "**"),
    @("post_ltx vs healthy", 1, 1736.16873682658, 11.0135363130856, 0.0207805893668765, 0.001, 0.0015, "**"),
    @("pre_ltx vs healthy , Country", 1, 738.576413458271, 4.80132620769844, 0.0113473740441028, 0.001, 0.001, "***"),
    @("pre_ltx vs post_ltx , Country", 1, 1532.67469153163, 9.8793465538545, 0.016102657510029, 0.001, 0.001, "***"),
    @("post_ltx vs healthy , Country", 1, 1728.7458235653, 10.9664484217891, 0.0206917428687702, 0.001, 0.001, "***"),
    @("pre_ltx vs healthy : Country", 1, 318.719016026873, 2.07739496027165, 0.00489674977960781, 0.325, 0.4875, "'"),
    @("pre_ltx vs post_ltx : Country", 1, 322.405182713765, 2.08193228651677, 0.00338726819551645, 0.524, 0.524, "'"),
    @("post_ltx vs healthy : Country", 1, 398.820084185802, 2.53760726113271, 0.00477356620064291, 0.112, 0.336, "'")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# The "sig" column is blank for rows 8-10 (no asterisks). Excel still records
# these as (empty) text cells rather than leaving them completely blank, so
# write them as an empty string and then strip the quote-prefix formatting
# that results from the leading apostrophe, restoring the default style.
foreach ($r in 8, 9, 10) {
    $cell = $ws.Cells.Item($r, 8)
    $cell.Style = "Normal"
}
